$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-57 down to 17-58
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with data (same as former row 16, but with an updated
# sampling date and a different "Primera" unit price)
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44600
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100112012
$ws.Cells.Item(16, 7).Value = "Espinaca"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 320
$ws.Cells.Item(16, 11).Value = 1400
$ws.Cells.Item(16, 12).Value = 1500
$ws.Cells.Item(16, 13).Value = 1450
$ws.Cells.Item(16, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(16, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 16).Value = 483
$ws.Cells.Item(16, 17).Value = 3
$ws.Cells.Item(16, 18).Value = "Hortaliza"
